$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "EAP - MEC-3A"
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = "EAP - MCT-2A"
$ws.Range("C6").Value = "EAP - MCT-2A"
$ws.Range("E6").Value = "-"
